# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" sheet (fund-level holdings) positioned right
# before the "总计" (Total) sheet, and prepends a matching summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# Use an existing quarter sheet as a style template for the new sheet.
$srcSheet = $wb.Worksheets.Item("2021-Q4")

# Create the new sheet positioned right before "总计".
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the outline properties used by the rest of the workbook's sheets
# (summary rows/columns below/right of the detail).
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Header row values
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy header formatting (bold/bordered/centered) from the template sheet.
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Force columns B, D, E, F, G to be stored as text (matching the rest
# of the workbook), so that leading zeros (e.g. fund code "012348") and
# exact decimal formatting (e.g. "38.10") are preserved instead of being
# coerced into numbers. The NumberFormat is reset back to "Normal" style
# afterwards so the cells don't end up with a stray custom style index.
$newSheet.Range("B2:B4").NumberFormat = "@"
$newSheet.Range("D2:G4").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "012348"
$newSheet.Range("C2").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）A"
$newSheet.Range("D2").Value = "38.10"
$newSheet.Range("E2").Value = "92.34"
$newSheet.Range("F2").Value = "6.81"
$newSheet.Range("G2").Value = "2.5946"
$newSheet.Range("H2").Value = 6

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "012349"
$newSheet.Range("C3").Value = "天弘恒生科技指数型发起式证券投资基金（QDII）C"
$newSheet.Range("D3").Value = "14.77"
$newSheet.Range("E3").Value = "92.34"
$newSheet.Range("F3").Value = "6.81"
$newSheet.Range("G3").Value = "1.0058"
$newSheet.Range("H3").Value = 6

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "539002"
$newSheet.Range("C4").Value = "建信新兴市场优选混合QDII"
$newSheet.Range("D4").Value = "0.14"
$newSheet.Range("E4").Value = "83.76"
$newSheet.Range("F4").Value = "4.24"
$newSheet.Range("G4").Value = "0.0059"
$newSheet.Range("H4").Value = 8

# Clean up the text-format styling so these data cells end up unstyled,
# matching the rest of the workbook's data rows.
$newSheet.Range("B2:B4").Style = "Normal"
$newSheet.Range("D2:G4").Style = "Normal"

# Copy the A-column style (centered/bordered index used for the leading
# numeric column) from the template sheet onto the new rows.
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# Now update the "总计" sheet: prepend a new summary row for 2022-Q1,
# shifting the existing rows down by one.
# Re-fetch the reference by name since the old COM reference can become
# stale/rebound after Worksheets.Add().
$totalSheet = $wb.Worksheets.Item("总计")

# Shift existing data rows (2:4) down to (3:5), carrying formatting along.
$totalSheet.Rows("2:4").Copy()
$totalSheet.Range("A3").PasteSpecial()

# The paste above does not reliably restore the formatting of the very
# last shifted row (A5), so fix it up explicitly by copying formats from
# the row right above it.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

# Write the new 2022-Q1 summary row (the cell styles already in place on
# row 2 match what is required, so only the values need to change).
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 3.61

# Renumber the index column (A) for the shifted rows 3..5.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the originally active sheet/tab (adding a sheet makes the new
# one active by default).
$wb.Worksheets.Item("2020-Q4").Activate()

$wb.Save()
